$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.012.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.727.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.28%  '

$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4847'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3470'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.21'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07233'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.050'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.00%  '

$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.870'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.746.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.812'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.74%  '

$ws.Range("E18").Value = '  -1.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06404'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.38%  '

$ws.Range("E20").Value = '  +0.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.709'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.075.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.073'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.922.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.063'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.035'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09318'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.635'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.374'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05903'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02176'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.430'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.1997'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.742'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5974'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.122'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.481'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.49%  '

$ws.Range("E46").Value = '  -4.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5612'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.846'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.100'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06652'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.31%  '
